# Update the scraped timestamp for every data row (2..399) from
# "2023-01-05 12:56:40" to "2023-01-05 20:49:40", plus a handful of
# ratingAmount (col D) bumps and productAriaLabel (col M) text tweaks
# that happened in this particular crawl run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 399; $r++) {
    $ws.Cells.Item($r, 15).Value = "2023-01-05 20:49:40"
}

# ratingAmount (column D) bumps
$ws.Cells.Item(5, 4).Value = 23
$ws.Cells.Item(17, 4).Value = 23
$ws.Cells.Item(21, 4).Value = 6
$ws.Cells.Item(26, 4).Value = 38
$ws.Cells.Item(96, 4).Value = 4
$ws.Cells.Item(114, 4).Value = 8
$ws.Cells.Item(129, 4).Value = 6
$ws.Cells.Item(214, 4).Value = 14
$ws.Cells.Item(266, 4).Value = 5

# productAriaLabel (column M) text tweaks
$ws.Cells.Item(32, 13).Value = "Roland Zwieback Classic 20% ab 2 Stück Aktion 4.40 Schweizer Franken"
$ws.Cells.Item(47, 13).Value = "Ölz Premium Buttertoast 2x  500g 50% Aktion 3.80 Schweizer Franken statt 7.60 Schweizer Franken"
$ws.Cells.Item(328, 13).Value = "Naturaplan Bio Haselnussrollen 4x60g 3 + 1 Aktion 3.05 Schweizer Franken statt 4.10 Schweizer Franken"
$ws.Cells.Item(337, 13).Value = "Pasquier Pancakes Choco 10 Stück 4.30 Schweizer Franken"
